$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Risks")

$ws.Range("A5").Value = "74d68dde-c792-4fdf-8cea-0a6960ef2e5c"
$ws.Range("B5").Value = 0.465
$ws.Range("C5").Value = 0.131
$ws.Range("D5").Value = 0.8
$ws.Range("E5").Value = "Mitigation needed"
